$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '25.976.10'
$ws.Range('E2').Value = '  +0.05%  '
$ws.Range('D3').Value = '1.633.32'
$ws.Range('E3').Value = '  -0.42%  '
$ws.Range('E4').Value = '  +0.27%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.00'
$ws.Range('E6').Value = '  -0.60%  '
$ws.Range('E7').Value = '  +0.27%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.251'
$ws.Range('E8').Value = '  -1.93%  '
$ws.Range('E9').Value = '  -2.38%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.53'
$ws.Range('E10').Value = '  -5.56%  '
$ws.Range('E11').Value = '  -0.33%  '
$ws.Range('B12').Value = 'WrappedEther'
$ws.Range('C12').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D12').Value = '1.681.04'
$ws.Range('E12').Value = '  +2.31%  '
$ws.Range('B13').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C13').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D13').Value = '1.860.21'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.20'
$ws.Range('E14').Value = '  -2.07%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.530'
$ws.Range('E15').Value = '  -2.75%  '
$ws.Range('D16').Value = '25.993.06'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = '0.0₃0745'
$ws.Range('E17').Value = '  -2.75%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.75'
$ws.Range('E18').Value = '  -1.85%  '
$ws.Range('E19').Value = '  +0.16%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '190.45'
$ws.Range('E20').Value = '  -1.39%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.25'
$ws.Range('E21').Value = '  -2.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.57'
$ws.Range('E22').Value = '  -3.68%  '
$ws.Range('E23').Value = '  -2.02%  '
$ws.Range('E24').Value = '  +0.62%  '
$ws.Range('E25').Value = '  -0.72%  '
$ws.Range('E26').Value = '  +0.34%  '
$ws.Range('E27').Value = '  -2.66%  '
$ws.Range('E28').Value = '  -2.28%  '
$ws.Range('E29').Value = '  -2.11%  '
$ws.Range('E30').Value = '  -1.15%  '
$ws.Range('E31').Value = '  -2.77%  '
$ws.Range('E32').Value = '  -3.12%  '
$ws.Range('E33').Value = '  -4.14%  '
$ws.Range('E34').Value = '  -1.59%  '
$ws.Range('E35').Value = '  -1.90%  '
$ws.Range('B36').Value = 'Maker'
$ws.Range('C36').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D36').Value = '1.135.41'
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('B37').Value = 'ARBITRUM'
$ws.Range('C37').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.868'
$ws.Range('E37').Value = '  -3.94%  '
$ws.Range('E38').Value = '  -1.49%  '
$ws.Range('E39').Value = '  -2.96%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0155'
$ws.Range('E40').Value = '  -1.39%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.64'
$ws.Range('E41').Value = '  -0.65%  '
$ws.Range('E42').Value = '  -2.00%  '
$ws.Range('E43').Value = '  -4.77%  '
$ws.Range('D44').Value = '1.770.76'
$ws.Range('E45').Value = '  -0.79%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.11'
$ws.Range('E46').Value = '  -2.71%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0526'
$ws.Range('E47').Value = '  -0.72%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.47'
$ws.Range('E48').Value = '  +1.23%  '
$ws.Range('E49').Value = '  -0.16%  '
$ws.Range('E50').Value = '  -2.50%  '
$ws.Range('E51').Value = '  +0.34%  '
